# Add data for 2022-03-27
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the running "through" date references.
$ws.Name = "Through 2022-03-19"
$ws.Range("B1").Value = "March 2022 (through March 19)"

# Update / add counts for column B (March 2022) and other month columns
# per neighborhood (row numbers match the sheet's existing rows).

# Austin (row 3)
$ws.Range("B3").Value = 7
$ws.Range("K3").Value = 5

# North Lawndale (row 4)
$ws.Range("B4").Value = 6

# Garfield Park (row 5)
$ws.Range("N5").Value = 5

# Rogers Park (row 6)
$ws.Range("T6").Value = 2

# Englewood (row 11)
$ws.Range("N11").Value = 5
$ws.Range("T11").Value = 4

# Lake View (row 12) - new value
$ws.Range("K12").Value = 1

# Washington Heights (row 18)
$ws.Range("H18").Value = 3

# Lincoln Park (row 19)
$ws.Range("B19").Value = 2

# Chatham (row 22)
$ws.Range("E22").Value = 2

# New City (row 32)
$ws.Range("H32").Value = 1
$ws.Range("T32").Value = 2

# Morgan Park (row 40) - new value
$ws.Range("H40").Value = 1

# South Chicago (row 84) - new value
$ws.Range("H84").Value = 1

# Streeterville (row 86) - new value
$ws.Range("Q86").Value = 1
